# cyd_hist.xlsx maintenance pass:
#  - some improvements to image floating (chart resized/repositioned)
#  - updated sizes in histogram spreadsheets (chart anchor shrunk)
#  - selection cursor left on a different cell when the file was saved

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Resize/reposition the histogram chart -------------------------------
# The chart's top-left anchor (column 5 / row 4, i.e. 0-indexed col 4 / row 3)
# stays put; only the bottom-right anchor moves from (col 14, row 23) to
# (col 11, row 15) in 1-indexed terms. Compute the new width/height in
# points from the actual column widths / row heights so the exported
# two-cell anchor offsets land on the target cell/offset exactly.
$co = $ws.ChartObjects().Item(1)

$fromColIdx = 5   # 1-indexed column the chart currently starts at
$fromRowIdx = 4   # 1-indexed row the chart currently starts at
$toColIdx   = 11  # 1-indexed column the chart should now end at
$toRowIdx   = 15  # 1-indexed row the chart should now end at

$fromLeft = 0.0
for ($c = 1; $c -lt $fromColIdx; $c++) {
    $fromLeft += $ws.Cells.Item(1, $c).Width
}
$fromLeft += 247650 / 12700.0   # original colOff (EMU -> pt)

$fromTop = 0.0
for ($r = 1; $r -lt $fromRowIdx; $r++) {
    $fromTop += $ws.Cells.Item($r, 1).Height
}
$fromTop += 104775 / 12700.0    # original rowOff (EMU -> pt)

$toLeft = 0.0
for ($c = 1; $c -lt $toColIdx; $c++) {
    $toLeft += $ws.Cells.Item(1, $c).Width
}
$toLeft += 190050 / 12700.0     # target colOff (EMU -> pt)

$toTop = 0.0
for ($r = 1; $r -lt $toRowIdx; $r++) {
    $toTop += $ws.Cells.Item($r, 1).Height
}
$toTop += 169275 / 12700.0      # target rowOff (EMU -> pt)

$co.Width = $toLeft - $fromLeft
$co.Height = $toTop - $fromTop

# --- Move the saved cell selection on the data sheet ----------------------
$ws.Range("S11").Select()
